$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = 2.84
$ws.Range("H3").Value = 2.66
$ws.Range("I3").Value = 2.72
$ws.Range("W3").Value = 1.54
$ws.Range("Z3").Value = 18.5
$ws.Range("AE3").Value = 29
$ws.Range("AJ3").Value = 46
$ws.Range("AL3").Value = 42
$ws.Range("F4").Value = 3.65
$ws.Range("G4").Value = 3.75
$ws.Range("H4").Value = 2.1
$ws.Range("I4").Value = 2.14
$ws.Range("AA4").Value = 26
$ws.Range("F5").Value = 1.58
$ws.Range("G5").Value = 1.75
$ws.Range("H5").Value = 2.22
$ws.Range("I5").Value = 7.4
$ws.Range("P5").Value = 2.16
$ws.Range("Q5").Value = 1.71
$ws.Range("F6").Value = 1.72
$ws.Range("G6").Value = 1.9
$ws.Range("H6").Value = 3.95
$ws.Range("I6").Value = 6.2
$ws.Range("K6").Value = 5.1
$ws.Range("P6").Value = 2
$ws.Range("Q6").Value = 1.78
$ws.Range("F7").Value = 1.77
$ws.Range("G7").Value = 1.95
$ws.Range("H7").Value = 4.1
$ws.Range("J7").Value = 3.85
$ws.Range("K7").Value = 4.6
$ws.Range("P7").Value = 2.22
$ws.Range("Q7").Value = 1.65
$ws.Range("F8").Value = 1.97
$ws.Range("G8").Value = 2.2
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 4.3
$ws.Range("J8").Value = 3.95
$ws.Range("K8").Value = 4.7
$ws.Range("P8").Value = 2.5
$ws.Range("Q8").Value = 1.54
$ws.Range("I9").Value = 5
$ws.Range("P9").Value = 1.68
$ws.Range("Q9").Value = 2.26
$ws.Range("F10").Value = 2.82
$ws.Range("G10").Value = 3.15
$ws.Range("I10").Value = 3
$ws.Range("P10").Value = 1.61
$ws.Range("Q10").Value = 2.38
$ws.Range("I11").Value = 4.6
$ws.Range("F12").Value = 2.92
$ws.Range("F13").Value = 2.84
$ws.Range("G13").Value = 3.1
$ws.Range("H13").Value = 2.72
$ws.Range("I13").Value = 2.88
$ws.Range("K13").Value = 3.45
$ws.Range("J20").Value = 7
$ws.Range("N20").Value = 6.2
$ws.Range("Z20").Value = 250
$ws.Range("AD20").Value = 75
$ws.Range("AG20").Value = 11.5
$ws.Range("AH20").Value = 40
$ws.Range("AI20").Value = 300
$ws.Range("I21").Value = 4.2
$ws.Range("AK21").Value = 18.5
$ws.Range("AM21").Value = 80
$ws.Range("F22").Value = 1.95
$ws.Range("G22").Value = 1.98
$ws.Range("H22").Value = 3.75
$ws.Range("W22").Value = 2.02
$ws.Range("X22").Value = 36
$ws.Range("AD22").Value = 17
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 7.8
$ws.Range("H23").Value = 1.49
$ws.Range("I23").Value = 1.5
$ws.Range("J23").Value = 5.1
$ws.Range("K23").Value = 5.3
$ws.Range("AC23").Value = 12.5
$ws.Range("AF23").Value = 1000
$ws.Range("F24").Value = 1.89
$ws.Range("I24").Value = 4.8
$ws.Range("J24").Value = 3.8
$ws.Range("G25").Value = 1.93
$ws.Range("H25").Value = 4.5
$ws.Range("I25").Value = 4.8
$ws.Range("K25").Value = 3.9
$ws.Range("P25").Value = 2.08
$ws.Range("Q25").Value = 1.87
$ws.Range("S25").Value = 3.15
$ws.Range("T25").Value = 1.79
$ws.Range("U25").Value = 2.18
$ws.Range("X25").Value = 15.5
$ws.Range("Y25").Value = 1000
$ws.Range("Z25").Value = 1000
$ws.Range("AA25").Value = 1000
$ws.Range("AD25").Value = 18.5
$ws.Range("AE25").Value = 1000
$ws.Range("AF25").Value = 12.5
$ws.Range("AH25").Value = 18.5
$ws.Range("AI25").Value = 1000
$ws.Range("AK25").Value = 19.5
$ws.Range("AM25").Value = 1000
$ws.Range("AO25").Value = 1000
$ws.Range("F26").Value = 2.4
$ws.Range("I26").Value = 3.45
$ws.Range("J26").Value = 3.45
$ws.Range("Y26").Value = 12.5
$ws.Range("AA26").Value = 1000
$ws.Range("AE26").Value = 40
$ws.Range("AF26").Value = 15.5
$ws.Range("AJ26").Value = 34
$ws.Range("AK26").Value = 27
$ws.Range("Q27").Value = 2.08
$ws.Range("J28").Value = 2.98
$ws.Range("F29").Value = 1.36
$ws.Range("G29").Value = 1.41
$ws.Range("H29").Value = 10.5
$ws.Range("I29").Value = 14
$ws.Range("J29").Value = 5
$ws.Range("K29").Value = 5.6
$ws.Range("P29").Value = 2
$ws.Range("Q29").Value = 1.86
$ws.Range("Q30").Value = 2.16
